# The workbook gained one new weekly data row. A new row was inserted at
# spreadsheet row 289 (pushing the previous rows 289-383 down to 290-384),
# and the new row was filled in with a fresh "Espinaca" price observation
# for "Terminal La Palmera de La Serena".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 289, shifting existing rows 289-383 down to 290-384
$ws.Rows.Item(289).Insert()

# Populate the newly inserted row 289 with the new observation's data
$ws.Range("A289").Value = 8
$ws.Range("B289").Value = "Terminal La Palmera de La Serena"
$ws.Range("C289").Value = "Coquimbo"
$ws.Range("D289").Value = 44988
$ws.Range("E289").Value = 4
$ws.Range("F289").Value = 100112012
$ws.Range("G289").Value = "Espinaca"
$ws.Range("H289").Value = "Sin especificar"
$ws.Range("I289").Value = "Primera"
$ws.Range("J289").Value = 1800
$ws.Range("K289").Value = 500
$ws.Range("L289").Value = 600
$ws.Range("M289").Value = 550
$ws.Range("N289").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O289").Value = "Provincia del Elquí"
$ws.Range("P289").Value = 1100
$ws.Range("Q289").Value = 0.5
$ws.Range("R289").Value = "Hortaliza"
